$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the cells we are about to write, then set each value explicitly
# so Excel does not reinterpret numeric-looking strings (e.g. "2.20" -> 2.2) or
# thousand-separated price strings as numbers/dates.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '43.528.02'

$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.73%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.419.29'

$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +2.04%  '

$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.02%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '306.62'

$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +1.05%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '97.38'

$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.56%  '

$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.31%  '

$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -0.03%  '

$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -1.91%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '34.96'

$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +2.18%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0797'

$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +2.39%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '18.51'

$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -0.27%  '

$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +1.70%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.788.56'

$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +1.97%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.404.36'

$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +2.40%  '

$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +3.42%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '43.511.16'

$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +0.72%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.43'

$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +2.23%  '

$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -2.15%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0₃0900'

$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +1.09%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '68.16'

$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -0.15%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '237.72'

$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +0.77%  '

$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +1.20%  '

$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +0.64%  '

$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.08%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '24.96'

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.20'

$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -6.99%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.43'

$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +3.00%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '32.32'

$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +2.48%  '

$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +17.93%  '

$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +0.62%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '18.39'

$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +6.14%  '

$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +0.03%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0752'

$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +2.91%  '

$ws.Range('B36').NumberFormat = '@'
$ws.Range('B36').Value = 'Monero'

$ws.Range('C36').NumberFormat = '@'
$ws.Range('C36').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '132.45'

$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +26.92%  '

$ws.Range('B37').NumberFormat = '@'
$ws.Range('B37').Value = 'ARBITRUM'

$ws.Range('C37').NumberFormat = '@'
$ws.Range('C37').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.90'

$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +2.79%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.93'

$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +5.16%  '

$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -0.39%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.27'

$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -1.31%  '

$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -0.10%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '20.95'

$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -6.52%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.947.70'

$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +0.10%  '

$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +0.80%  '

$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +2.39%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.84'

$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +3.11%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.31'

$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -1.67%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.647.19'

$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +2.12%  '

$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +3.00%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '52.67'

$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -0.76%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '72.36'

$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +0.04%  '
